# Updates the "StructureDefinition-employee-pay-frequency" workbook:
#   - Metadata sheet: URL / Version / Date / Publisher values (IBM -> LinuxForHealth rebrand)
#   - Elements sheet: the "Extension" row's Constraint(s) cell is cleared; the
#     combined ele-1/ext-1 constraint now only shows on the "Extension.extension" row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-pay-frequency"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" element; its Constraint(s) value (column AI) is cleared.
$elements.Range("AI2").Value = ""

# Row 4 = "Extension.extension" element; its Constraint(s) value (column AI)
# now carries the ele-1/ext-1 constraint text that used to be shown on row 2.
$elements.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
